$d = $word.ActiveDocument

# Locate the unique spot right before "Date:" in the "Ref: MSL/..." paragraph -
# that text is preceded by a run holding a tab + 6 spaces.
$ins = $d.Content
$found = $ins.Find.Execute('Date:        .${month}.${year}', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Date:' anchor text"
}
$ins.Collapse(1)

# Insert three additional runs (tab / tab / tab + 8 spaces), all using the
# same sz=22 / szCs=22 direct character formatting as the surrounding runs.
# A transient bookmark is dropped after each insertion so the engine keeps
# each insertion in its own run instead of folding it back into its
# same-format neighbour.
$ins.InsertAfter([char]9)
$d.Bookmarks.Add("zzTmpRunBreak1", $ins)
$ins.Collapse(0)

$ins.InsertAfter([char]9)
$d.Bookmarks.Add("zzTmpRunBreak2", $ins)
$ins.Collapse(0)

$ins.InsertAfter([char]9 + "        ")
$d.Bookmarks.Add("zzTmpRunBreak3", $ins)
$ins.Collapse(0)

$d.Bookmarks("zzTmpRunBreak1").Delete()
$d.Bookmarks("zzTmpRunBreak2").Delete()
$d.Bookmarks("zzTmpRunBreak3").Delete()
